# Payslip process / Config / Config.xlsx
# Fill in the three new settings rows (payroll template, save location,
# DKK conversion URL) on the "Settings" sheet and move the active
# selection to where the user left off (C17).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Row 9 - Payslip payroll template
$ws.Range("A9").Value = "Payslip payroll template"
$ws.Range("B9").Value = 'C:\\Users\\premal\\Google Drive\\QAC projects\\HR automation\\automation\\HR-Payroll-automation\\Payslip process\\DataBases\\payslip_template_employees.xlsx'

# Row 10 - Payslip save location
$ws.Range("A10").Value = "Payslip save location"
$ws.Range("B10").Value = 'C:\\Users\\premal\\Google Drive\\QAC projects\\HR automation\\automation\\HR-Payroll-automation\\Payslip process\\Payslips\\'

# Row 11 - URL used to scrape the GBP -> DKK conversion rate
$ws.Range("A11").Value = "URL DKK"
$ws.Range("B11").Value = 'https://www.xe.com/currencyconverter/convert/?Amount=1&From=GBP&To=DKK'
$ws.Range("C11").Value = "Website to scrape DKK conversion rate"

# Leave the selection where the author left it when they saved.
$ws.Range("C17").Select()
